# Update profit-calculation cells (currentAveragePrice/.../LeveProfit columns)
# across the Anima_Profits sheets, per scheduled-runner recomputation.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 10043.889
$ws.Range("I62").Value = 3499.1667
$ws.Range("J62").Value = 23133.334
$ws.Range("K62").Value = 3499.1667
$ws.Range("L62").Value = 23133.334
$ws.Range("M62").Value = -2875.1667
$ws.Range("N62").Value = -24381.334
$ws.Range("H65").Value = 10043.889
$ws.Range("I65").Value = 3499.1667
$ws.Range("J65").Value = 23133.334
$ws.Range("K65").Value = 17495.8335
$ws.Range("L65").Value = 115666.67
$ws.Range("M65").Value = -14375.8335
$ws.Range("N65").Value = -121906.67
$ws.Range("H70").Value = 1628.4286
$ws.Range("I70").Value = 1499.75
$ws.Range("J70").Value = 1800
$ws.Range("K70").Value = 4499.25
$ws.Range("L70").Value = 5400
$ws.Range("M70").Value = -4229.25
$ws.Range("N70").Value = -5940
$ws.Range("H73").Value = 1628.4286
$ws.Range("I73").Value = 1499.75
$ws.Range("J73").Value = 1800
$ws.Range("K73").Value = 4499.25
$ws.Range("L73").Value = 5400
$ws.Range("M73").Value = -3563.25
$ws.Range("N73").Value = -7272
$ws.Range("H98").Value = 1187.7858
$ws.Range("I98").Value = 902.5833
$ws.Range("J98").Value = 2899
$ws.Range("K98").Value = 902.5833
$ws.Range("L98").Value = 2899
$ws.Range("M98").Value = 595.4167
$ws.Range("N98").Value = -5895
$ws.Range("H113").Value = 2479.7368
$ws.Range("I113").Value = 2216.111
$ws.Range("J113").Value = 2717
$ws.Range("K113").Value = 2216.111
$ws.Range("L113").Value = 2717
$ws.Range("M113").Value = 1037.889
$ws.Range("N113").Value = -9225
$ws.Range("H122").Value = 1187.7858
$ws.Range("I122").Value = 902.5833
$ws.Range("J122").Value = 2899
$ws.Range("K122").Value = 2707.7499
$ws.Range("L122").Value = 8697
$ws.Range("M122").Value = -257.7498999999998
$ws.Range("N122").Value = -13597
$ws.Range("H132").Value = 2552.0244
$ws.Range("I132").Value = 2368.3242
$ws.Range("J132").Value = 4251.25
$ws.Range("K132").Value = 7104.9726
$ws.Range("L132").Value = 12753.75
$ws.Range("M132").Value = -4574.9726
$ws.Range("N132").Value = -17813.75
$ws.Range("H135").Value = 1849
$ws.Range("I135").Value = 989.0909
$ws.Range("K135").Value = 8901.8181
$ws.Range("M135").Value = -6366.8181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1141.5834
$ws.Range("I97").Value = 1229
$ws.Range("J97").Value = 180
$ws.Range("K97").Value = 1229
$ws.Range("L97").Value = 180
$ws.Range("M97").Value = -733
$ws.Range("N97").Value = -1172
$ws.Range("H122").Value = 46930.91
$ws.Range("I122").Value = 57030.555
$ws.Range("J122").Value = 1482.5
$ws.Range("K122").Value = 171091.665
$ws.Range("L122").Value = 4447.5
$ws.Range("M122").Value = -168641.665
$ws.Range("N122").Value = -9347.5
$ws.Range("H132").Value = 3769.889
$ws.Range("I132").Value = 3401.7437
$ws.Range("J132").Value = 4727.067
$ws.Range("K132").Value = 10205.2311
$ws.Range("L132").Value = 14181.201
$ws.Range("M132").Value = -7675.231100000001
$ws.Range("N132").Value = -19241.201

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2490.2173
$ws.Range("I86").Value = 2185.2
$ws.Range("J86").Value = 3062.125
$ws.Range("K86").Value = 2185.2
$ws.Range("L86").Value = 3062.125
$ws.Range("M86").Value = -1062.2
$ws.Range("N86").Value = -5308.125
$ws.Range("H89").Value = 2490.2173
$ws.Range("I89").Value = 2185.2
$ws.Range("J89").Value = 3062.125
$ws.Range("K89").Value = 10926
$ws.Range("L89").Value = 15310.625
$ws.Range("M89").Value = -5310
$ws.Range("N89").Value = -26542.625
$ws.Range("H134").Value = 3744.8215
$ws.Range("I134").Value = 3959.8096
$ws.Range("J134").Value = 3099.8572
$ws.Range("K134").Value = 11879.4288
$ws.Range("L134").Value = 9299.571599999999
$ws.Range("M134").Value = -9344.4288
$ws.Range("N134").Value = -14369.5716

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4933.2437
$ws.Range("I31").Value = 1337.2572
$ws.Range("J31").Value = 7860.2095
$ws.Range("K31").Value = 1337.2572
$ws.Range("L31").Value = 7860.2095
$ws.Range("M31").Value = -1042.2572
$ws.Range("N31").Value = -8450.209500000001
$ws.Range("H34").Value = 4933.2437
$ws.Range("I34").Value = 1337.2572
$ws.Range("J34").Value = 7860.2095
$ws.Range("K34").Value = 1337.2572
$ws.Range("L34").Value = 7860.2095
$ws.Range("M34").Value = -1135.2572
$ws.Range("N34").Value = -8264.209500000001
$ws.Range("H94").Value = 125001576
$ws.Range("I94").Value = 1000000000
$ws.Range("J94").Value = 1804
$ws.Range("K94").Value = 1000000000
$ws.Range("L94").Value = 1804
$ws.Range("M94").Value = -999999549
$ws.Range("N94").Value = -2706
$ws.Range("H99").Value = 2059.9333
$ws.Range("I99").Value = 1652.875
$ws.Range("J99").Value = 2147.946
$ws.Range("K99").Value = 1652.875
$ws.Range("L99").Value = 2147.946
$ws.Range("M99").Value = -154.875
$ws.Range("N99").Value = -5143.946
$ws.Range("H122").Value = 1820.2084
$ws.Range("I122").Value = 1300.5555
$ws.Range("J122").Value = 2132
$ws.Range("K122").Value = 3901.6665
$ws.Range("L122").Value = 6396
$ws.Range("M122").Value = -1451.6665
$ws.Range("N122").Value = -11296
$ws.Range("H126").Value = 2059.9333
$ws.Range("I126").Value = 1652.875
$ws.Range("J126").Value = 2147.946
$ws.Range("K126").Value = 4958.625
$ws.Range("L126").Value = 6443.838
$ws.Range("M126").Value = -2488.625
$ws.Range("N126").Value = -11383.838

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 870.3889
$ws.Range("I113").Value = 802.8333
$ws.Range("J113").Value = 1005.5
$ws.Range("K113").Value = 2408.4999
$ws.Range("L113").Value = 3016.5
$ws.Range("M113").Value = -238.4998999999998
$ws.Range("N113").Value = -7356.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6453.2856
$ws.Range("I122").Value = 9277.23
$ws.Range("J122").Value = 1864.375
$ws.Range("K122").Value = 27831.69
$ws.Range("L122").Value = 5593.125
$ws.Range("M122").Value = -25381.69
$ws.Range("N122").Value = -10493.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7288.9
$ws.Range("I7").Value = 7209.8887
$ws.Range("J7").Value = 8000
$ws.Range("K7").Value = 7209.8887
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = -7097.8887
$ws.Range("N7").Value = -8224
$ws.Range("H16").Value = 797.1177
$ws.Range("I16").Value = 797.1177
$ws.Range("K16").Value = 797.1177
$ws.Range("M16").Value = -627.1177
$ws.Range("H40").Value = 2628.9167
$ws.Range("I40").Value = 2364.7
$ws.Range("K40").Value = 2364.7
$ws.Range("M40").Value = -2228.7
$ws.Range("H46").Value = 592.46155
$ws.Range("I46").Value = 540.2
$ws.Range("J46").Value = 625.125
$ws.Range("K46").Value = 540.2
$ws.Range("L46").Value = 625.125
$ws.Range("M46").Value = -352.2
$ws.Range("N46").Value = -1001.125
$ws.Range("H93").Value = 6849.9473
$ws.Range("I93").Value = 7571.8125
$ws.Range("J93").Value = 3000
$ws.Range("K93").Value = 7571.8125
$ws.Range("L93").Value = 3000
$ws.Range("M93").Value = -6323.8125
$ws.Range("N93").Value = -5496
$ws.Range("H122").Value = 4638.8696
$ws.Range("I122").Value = 3734
$ws.Range("K122").Value = 11202
$ws.Range("M122").Value = -8752
$ws.Range("H126").Value = 7288.9
$ws.Range("I126").Value = 7209.8887
$ws.Range("J126").Value = 8000
$ws.Range("K126").Value = 21629.6661
$ws.Range("L126").Value = 24000
$ws.Range("M126").Value = -19159.6661
$ws.Range("N126").Value = -28940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1880.4615
$ws.Range("I122").Value = 1896.6086
$ws.Range("J122").Value = 1756.6666
$ws.Range("K122").Value = 5689.825800000001
$ws.Range("L122").Value = 5269.9998
$ws.Range("M122").Value = -3239.825800000001
$ws.Range("N122").Value = -10169.9998
$ws.Range("H126").Value = 1251.5333
$ws.Range("I126").Value = 1269.5714
$ws.Range("J126").Value = 999
$ws.Range("K126").Value = 3808.7142
$ws.Range("L126").Value = 2997
$ws.Range("M126").Value = -1338.7142
$ws.Range("N126").Value = -7937
$ws.Range("H132").Value = 4507679
$ws.Range("I132").Value = 3449.5293
$ws.Range("K132").Value = 10348.5879
$ws.Range("M132").Value = -7818.5879
$ws.Range("H136").Value = 2858.7358
$ws.Range("I136").Value = 2459.389
$ws.Range("J136").Value = 3704.4119
$ws.Range("K136").Value = 7378.167
$ws.Range("L136").Value = 11113.2357
$ws.Range("M136").Value = -4828.167
$ws.Range("N136").Value = -16213.2357
